$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writes, then re-protect afterward.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (cell A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."
$ws.Rows.Item(59).AutoFit()

# Update Weight (col D) and Percent Change (col E) values for rows 2-56.
$ws.Cells.Item(2, 4).Value = 0.02442776705258207
$ws.Cells.Item(2, 5).Value = -0.02937946183415707
$ws.Cells.Item(3, 4).Value = 0.01809919892561518
$ws.Cells.Item(3, 5).Value = -0.02096436058700213
$ws.Cells.Item(4, 4).Value = 0.01809804911309211
$ws.Cells.Item(4, 5).Value = -0.02465057179161367
$ws.Cells.Item(5, 4).Value = 0.0202756023963626
$ws.Cells.Item(5, 5).Value = -0.01640785232932895
$ws.Cells.Item(6, 4).Value = 0.01952381664169657
$ws.Cells.Item(6, 5).Value = -0.01813898704358086
$ws.Cells.Item(7, 4).Value = 0.02677760057889227
$ws.Cells.Item(7, 5).Value = -0.005524861878452914
$ws.Cells.Item(8, 4).Value = 0.0188254971693532
$ws.Cells.Item(8, 5).Value = 0.005008347245408995
$ws.Cells.Item(9, 4).Value = 0.01976604381322292
$ws.Cells.Item(9, 5).Value = -0.02039866594275974
$ws.Cells.Item(10, 4).Value = 0.01932144963763658
$ws.Cells.Item(10, 5).Value = -0.02915972387526788
$ws.Cells.Item(11, 4).Value = 0.01969053945754145
$ws.Cells.Item(11, 5).Value = -0.02189781021897808
$ws.Cells.Item(12, 4).Value = 0.01928235601185226
$ws.Cells.Item(12, 5).Value = -0.01282051282051266
$ws.Cells.Item(13, 4).Value = 0.01950771926637362
$ws.Cells.Item(13, 5).Value = -0.005344021376085584
$ws.Cells.Item(14, 4).Value = 0.01861661456099583
$ws.Cells.Item(14, 5).Value = -0.01052024787433348
$ws.Cells.Item(15, 4).Value = 0.01697353246553149
$ws.Cells.Item(15, 5).Value = 0.003974167908594284
$ws.Cells.Item(16, 4).Value = 0.01764310662479816
$ws.Cells.Item(16, 5).Value = -0.006212934199378628
$ws.Cells.Item(17, 4).Value = 0.01586262193282718
$ws.Cells.Item(17, 5).Value = 0.001328903654485014
$ws.Cells.Item(18, 4).Value = 0.01540902089247681
$ws.Cells.Item(18, 5).Value = -0.04290617848970246
$ws.Cells.Item(19, 4).Value = 0.01656879845741151
$ws.Cells.Item(19, 5).Value = -0.01895674300254435
$ws.Cells.Item(20, 4).Value = 0.01845832370365345
$ws.Cells.Item(20, 5).Value = -0.01349667774086383
$ws.Cells.Item(21, 4).Value = 0.02000673790138518
$ws.Cells.Item(21, 5).Value = -0.03465517241379301
$ws.Cells.Item(22, 4).Value = 0.02073955174942059
$ws.Cells.Item(22, 5).Value = -0.001644736842105199
$ws.Cells.Item(23, 4).Value = 0.01982851696030962
$ws.Cells.Item(23, 5).Value = -0.03754711510582776
$ws.Cells.Item(24, 4).Value = 0.02143633813839986
$ws.Cells.Item(24, 5).Value = -0.03316645807259067
$ws.Cells.Item(25, 4).Value = 0.02120867525883238
$ws.Cells.Item(25, 5).Value = -0.02746855573225393
$ws.Cells.Item(26, 4).Value = 0.01966639339455702
$ws.Cells.Item(26, 5).Value = -0.026368101028999
$ws.Cells.Item(27, 4).Value = 0.02054829959975026
$ws.Cells.Item(27, 5).Value = -0.03517803517803519
$ws.Cells.Item(28, 4).Value = 0.02675671231805654
$ws.Cells.Item(28, 5).Value = 0.003545261167572811
$ws.Cells.Item(29, 4).Value = 0.01992625102477041
$ws.Cells.Item(29, 5).Value = -0.03808424697057122
$ws.Cells.Item(30, 4).Value = 0.01244384603090466
$ws.Cells.Item(30, 5).Value = 0.004404404404404483
$ws.Cells.Item(31, 4).Value = 0.009021237420571909
$ws.Cells.Item(31, 5).Value = -0.01486988847583648
$ws.Cells.Item(32, 4).Value = 0.01610044148968177
$ws.Cells.Item(32, 5).Value = 0.0222457627118644
$ws.Cells.Item(33, 4).Value = 0.02044021722258186
$ws.Cells.Item(33, 5).Value = -0.02786840674279512
$ws.Cells.Item(34, 4).Value = 0.01841060648394613
$ws.Cells.Item(34, 5).Value = -0.03543212832176212
$ws.Cells.Item(35, 4).Value = 0.01810916396748177
$ws.Cells.Item(35, 5).Value = -0.009397024275646038
$ws.Cells.Item(36, 4).Value = 0.01642008937109471
$ws.Cells.Item(36, 5).Value = 0.04269175108538348
$ws.Cells.Item(37, 4).Value = 0.01950771926637362
$ws.Cells.Item(37, 5).Value = -0.0145290581162324
$ws.Cells.Item(38, 4).Value = 0.0197202429477207
$ws.Cells.Item(38, 5).Value = -0.01175841795831112
$ws.Cells.Item(39, 4).Value = 0.02530200784095486
$ws.Cells.Item(39, 5).Value = -0.01345128453708189
$ws.Cells.Item(40, 4).Value = 0.01760822897826509
$ws.Cells.Item(40, 5).Value = -0.003308519437551616
$ws.Cells.Item(41, 4).Value = 0.0228245451245841
$ws.Cells.Item(41, 5).Value = -0.0212923159591617
$ws.Cells.Item(42, 4).Value = 0.01958533161168071
$ws.Cells.Item(42, 5).Value = -0.02361033649377198
$ws.Cells.Item(43, 4).Value = 0.01983771546049417
$ws.Cells.Item(43, 5).Value = -0.006172839506172978
$ws.Cells.Item(44, 4).Value = 0.01859572630016009
$ws.Cells.Item(44, 5).Value = -0.0213320692107134
$ws.Cells.Item(45, 4).Value = 0.02041875405548458
$ws.Cells.Item(45, 5).Value = -0.02158610980760212
$ws.Cells.Item(46, 4).Value = 0.01956214372579884
$ws.Cells.Item(46, 5).Value = -0.008620689655172376
$ws.Cells.Item(47, 4).Value = 0.01783589185783258
$ws.Cells.Item(47, 5).Value = -0.01489169675090263
$ws.Cells.Item(48, 4).Value = 0.01595939782018542
$ws.Cells.Item(48, 5).Value = 0.004322766570605152
$ws.Cells.Item(49, 4).Value = 0.0170310230916849
$ws.Cells.Item(49, 5).Value = -0.0007763975155280489
$ws.Cells.Item(50, 4).Value = 0.01749554735100442
$ws.Cells.Item(50, 5).Value = -0.01498422712933756
$ws.Cells.Item(51, 4).Value = 0.01639766802689488
$ws.Cells.Item(51, 5).Value = -0.02827024436990899
$ws.Cells.Item(52, 4).Value = 0.01804822390375916
$ws.Cells.Item(52, 5).Value = -0.01571458908473145
$ws.Cells.Item(53, 4).Value = 0.01555274745786033
$ws.Cells.Item(53, 5).Value = -0.008378718056137324
$ws.Cells.Item(54, 4).Value = 0.007376047335481949
$ws.Cells.Item(54, 5).Value = 0.003507404520654678
$ws.Cells.Item(55, 4).Value = 0.007150300810119576
$ws.Cells.Item(55, 5).Value = 0.003001715265866123
$ws.Cells.Item(56, 4).Value = 0.9999999999999999
$ws.Cells.Item(56, 5).Value = -0.01485049945939654

# Restore sheet protection.
$ws.Protect()
